$d = $word.ActiveDocument

function Replace-ParagraphXml([int]$index, [string]$bodyXml) {
    $p = $d.Paragraphs($index)
    $r = $p.Range
    $target = $d.Range($r.Start, $r.End - 1)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($pkg)
}

Replace-ParagraphXml 9 '<w:p w14:paraId="50DC3710" w14:textId="6F072BB4" w:rsidR="00EA21F0" w:rsidRDefault="00CC7E94" w:rsidP="00EA21F0"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">In the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SceneView</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> constructor, I changed the tilt element of the camera element from 60 to 75 to set the viewing angle to 75 degrees, z element of the camera&#8217;s position element from 290 to 250 to lower the camera&#8217;s height by </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>40, and</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> changed the heading element of the camera from 330 to 270 to point the camera in a different direction.</w:t></w:r></w:p>'
Replace-ParagraphXml 10 '<w:p w14:paraId="0124996C" w14:textId="0174B9B8" w:rsidR="005E6103" w:rsidRDefault="005E6103" w:rsidP="00EA21F0"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">I navigated to </w:t></w:r><w:hyperlink r:id="rId6" w:history="1"><w:r w:rsidRPr="004743CE"><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr><w:t>https://developers.arcgis.com/javascript/latest/visualization/symbols-color-ramps/esri-web-style-symbols-3d/</w:t></w:r></w:hyperlink><w:r><w:t xml:space="preserve"> and found a bunch of funny-named detailed tree objects</w:t></w:r><w:r w:rsidR="007B4415"><w:t xml:space="preserve">, then in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>uniqueValueInfos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> list, went through each element and changed their&#8217; symbol&#8217;s name element to the name of the funny tree object.</w:t></w:r></w:p>'
Replace-ParagraphXml 11 '<w:p w14:paraId="76056DFB" w14:textId="17AF43E7" w:rsidR="00A3173F" w:rsidRDefault="00A3173F" w:rsidP="00EA21F0"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">I changed the lighting of the environment element in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SceneView</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to be today at noon</w:t></w:r><w:r w:rsidR="00C4270A"><w:t>.</w:t></w:r></w:p>'
Replace-ParagraphXml 12 '<w:p w14:paraId="39CCA40C" w14:textId="4B13F271" w:rsidR="00C4270A" w:rsidRPr="00EA21F0" w:rsidRDefault="00C4270A" w:rsidP="00EA21F0"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">I created a new variable called </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>myLayer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:t xml:space="preserve">selected the URL for the Bike paths </w:t></w:r><w:r><w:t xml:space="preserve">from </w:t></w:r><w:r><w:t>the University of Iowa&#8217;s ArcGIS REST API service</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r w:rsidR="00E82D18"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00000EC6"><w:t xml:space="preserve">I then changed the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>map.add</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">() function to add feature layers into the map to be </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>map.addMany</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(), which takes many elements in a list and adds them to the map</w:t></w:r><w:r><w:t xml:space="preserve"> (I was a little confused because some of the trees I specified didn&#8217;t load, though I was using the correct names for each of them). </w:t></w:r></w:p>'
